# Updated cryptos list with GitHub Actions
# Applies per-cell value changes (Price and Volume(1h) columns),
# plus two row-content swaps (Filecoin/Stellar, Algorand/EnergySwap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.658.22"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.84"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.06"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5371"
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3200"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07064"
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.09"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7771"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07828"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.41"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.51"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.056"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.18"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008026"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.682.24"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.092.01"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.650"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.053"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.414"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.82"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.218"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.699"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.83"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.299"
$ws.Range("E30").Value = "  +2.85%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08757"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.122"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04886"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7397"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.894"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.356"
$ws.Range("E38").Value = "  +6.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01755"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4847"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9116"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.57"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.932"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.745"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4215"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.193"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1256"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.13"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05839"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8990"
$ws.Range("E51").Value = "  +0.71%  "
